# Generate Report for Handback
# Updates the "generate date" / handoff / handback timestamp cells that get
# refreshed whenever the handback status report is regenerated.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview!G2 - "Latest HO Xliff Generate Date" (shared with de-de!H2)
$wsOverview.Range("G2").Value = "2016-09-03 17:13:36"

# zh-cn!H2 - "Correspond Handoff Datetime"
$wsZhCn.Range("H2").Value = "2016-09-03 17:13:31"

# zh-cn!K2 - "Correspond Handback DateTime"
$wsZhCn.Range("K2").Value = "2016-09-03 17:13:50"

# de-de!H2 - "Correspond Handoff Datetime" (same shared string as Overview!G2)
$wsDeDe.Range("H2").Value = "2016-09-03 17:13:36"

# de-de!K2 - "Correspond Handback DateTime"
$wsDeDe.Range("K2").Value = "2016-09-03 17:13:57"

$wb.Save()
